# Typography / date fixes for "Sitting in meetings all day long.pptx"
#
# 1) Refresh the cached "datetimeFigureOut" field text (4/10/24 -> 5/20/24)
#    on the slide master, every custom (slide) layout, and the notes master.
# 2) Update the "Rev. 1 (2024-04-??)" stamp on the title slide to "2024-05-??".
# 3) Straighten curly quotation marks on a few slides.

$p = $ppt.ActivePresentation

function Set-DateText {
    param($shape, [string]$newText)
    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -ne $newText) {
        $len = $tr.Text.Length
        $sub = $tr.Characters(1, $len)
        $sub.Text = $newText
    }
}

# --- 1a. Slide master "Date Placeholder" ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        Set-DateText $sh "5/20/24"
    }
}

# --- 1b. Every custom (slide) layout's "Date Placeholder" ---
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $cl = $layouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            Set-DateText $sh "5/20/24"
        }
    }
}

# --- 1c. Notes master "Date Placeholder" ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        Set-DateText $sh "5/20/24"
    }
}

# --- 2. "Rev. 1 (2024-04-??), " -> "Rev. 1 (2024-05-??), " on slide 1 ---
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    $tr = $sh.TextFrame.TextRange
    if ($tr.Text -like "Rev. 1 (2024-04-??), *") {
        $prefix = "Rev. 1 (2024-04-??), "
        $sub = $tr.Characters(1, $prefix.Length)
        $sub.Text = "Rev. 1 (2024-05-??), "
    }
}

# --- 3a. Slide 19 title: curly double quotes -> straight quotes ---
$s19 = $p.Slides.Item(19)
for ($i = 1; $i -le $s19.Shapes.Count; $i++) {
    $sh = $s19.Shapes.Item($i)
    $tr = $sh.TextFrame.TextRange
    if ($tr.Text -like '*"Managing" in all sorts of directions*') {
        $target = '"Managing" in all sorts of directions'
        $full = $tr.Text
        $startPos = $full.IndexOf($target) + 1
        $sub = $tr.Characters($startPos, $target.Length)
        $sub.Text = $target
    }
}

# --- 3b. "Dealing with "large" amounts" -> straighten the opening quote ---
foreach ($slideIdx in 23, 24, 25) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        $tr = $sh.TextFrame.TextRange
        $target = 'Dealing with "large" amounts'
        if ($tr.Text.StartsWith($target)) {
            $sub = $tr.Characters(1, $target.Length)
            $sub.Text = $target
        }
    }
}
